$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-46 down to 43-47
$ws.Rows.Item(42).Insert()

# Match the style of the new A42 label cell to the other column-A bucket cells (bold, bordered, centered)
$ws.Range("A42").Font.Bold = $true
$ws.Range("A42").HorizontalAlignment = -4108
$ws.Range("A42").VerticalAlignment = -4160
$ws.Range("A42").Borders.LineStyle = 1

# Row 42
$ws.Cells.Item(42, 1).Value = 60
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 1.15
$ws.Cells.Item(42, 4).Value = 0.02
$ws.Cells.Item(42, 5).Value = 7.09
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 1.14
$ws.Cells.Item(42, 9).Value = 14.86
$ws.Cells.Item(42, 10).Value = 0.03
$ws.Cells.Item(42, 11).Value = 5.99
$ws.Cells.Item(42, 12).Value = 1.12
$ws.Cells.Item(42, 13).Value = 5.98
$ws.Cells.Item(42, 14).Value = 6
$ws.Cells.Item(42, 15).Value = 14.16
$ws.Cells.Item(42, 16).Value = 14.16
$ws.Cells.Item(42, 17).Value = 14.16
$ws.Cells.Item(42, 18).Value = 0
$ws.Cells.Item(42, 19).Value = 0
$ws.Cells.Item(42, 20).Value = 0
$ws.Cells.Item(42, 21).Value = 14.16
$ws.Cells.Item(42, 22).Value = 6.97
$ws.Cells.Item(42, 23).Value = 0
$ws.Cells.Item(42, 24).Value = 0
$ws.Cells.Item(42, 25).Value = 0
$ws.Cells.Item(42, 26).Value = 0
$ws.Cells.Item(42, 27).Value = 36.99
$ws.Cells.Item(42, 28).Value = 302.58
$ws.Cells.Item(42, 29).Value = 0.01
$ws.Cells.Item(42, 30).Value = 0
$ws.Cells.Item(42, 31).Value = 0
$ws.Cells.Item(42, 32).Value = 1.2
$ws.Cells.Item(42, 33).Value = 0.03
$ws.Cells.Item(42, 34).Value = 6.94
$ws.Cells.Item(42, 35).Value = 0
$ws.Cells.Item(42, 36).Value = 0
$ws.Cells.Item(42, 37).Value = 1.2
$ws.Cells.Item(42, 38).Value = 14.82
$ws.Cells.Item(42, 39).Value = 0.03
$ws.Cells.Item(42, 40).Value = 5.96
$ws.Cells.Item(42, 41).Value = 1.05
$ws.Cells.Item(42, 42).Value = 5.96
$ws.Cells.Item(42, 43).Value = 5.96
$ws.Cells.Item(42, 44).Value = 13.94
$ws.Cells.Item(42, 45).Value = 13.94
$ws.Cells.Item(42, 46).Value = 13.94
$ws.Cells.Item(42, 47).Value = 0
$ws.Cells.Item(42, 48).Value = 0
$ws.Cells.Item(42, 49).Value = 0
$ws.Cells.Item(42, 50).Value = 13.94
$ws.Cells.Item(42, 51).Value = 6.95
$ws.Cells.Item(42, 52).Value = 0.01
$ws.Cells.Item(42, 53).Value = 0
$ws.Cells.Item(42, 54).Value = 0
$ws.Cells.Item(42, 55).Value = 0.01
$ws.Cells.Item(42, 56).Value = 37.02
$ws.Cells.Item(42, 57).Value = 303.09
$ws.Cells.Item(42, 58).Value = 0.91
$ws.Cells.Item(42, 59).Value = 0

# Row 43
$ws.Cells.Item(43, 1).Value = 62
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(43, 3).Value = 0.42
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 0.42
$ws.Cells.Item(43, 9).Value = 17
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 0
$ws.Cells.Item(43, 14).Value = 0
$ws.Cells.Item(43, 15).Value = 17
$ws.Cells.Item(43, 16).Value = 17
$ws.Cells.Item(43, 17).Value = 17
$ws.Cells.Item(43, 18).Value = 0
$ws.Cells.Item(43, 19).Value = 0
$ws.Cells.Item(43, 20).Value = 0
$ws.Cells.Item(43, 21).Value = 17
$ws.Cells.Item(43, 22).Value = 0
$ws.Cells.Item(43, 23).Value = 0
$ws.Cells.Item(43, 24).Value = 0
$ws.Cells.Item(43, 25).Value = 0
$ws.Cells.Item(43, 26).Value = 0
$ws.Cells.Item(43, 27).Value = 42.99
$ws.Cells.Item(43, 28).Value = 322.17
$ws.Cells.Item(43, 29).Value = 0.55
$ws.Cells.Item(43, 30).Value = 0
$ws.Cells.Item(43, 31).Value = 0
$ws.Cells.Item(43, 32).Value = 0.3
$ws.Cells.Item(43, 33).Value = 0
$ws.Cells.Item(43, 34).Value = 0
$ws.Cells.Item(43, 35).Value = 0
$ws.Cells.Item(43, 36).Value = 0
$ws.Cells.Item(43, 37).Value = 0.3
$ws.Cells.Item(43, 38).Value = 17
$ws.Cells.Item(43, 39).Value = 0
$ws.Cells.Item(43, 40).Value = 0
$ws.Cells.Item(43, 41).Value = 0
$ws.Cells.Item(43, 42).Value = 0
$ws.Cells.Item(43, 43).Value = 0
$ws.Cells.Item(43, 44).Value = 17
$ws.Cells.Item(43, 45).Value = 17
$ws.Cells.Item(43, 46).Value = 17
$ws.Cells.Item(43, 47).Value = 0
$ws.Cells.Item(43, 48).Value = 0
$ws.Cells.Item(43, 49).Value = 0
$ws.Cells.Item(43, 50).Value = 17
$ws.Cells.Item(43, 51).Value = 0
$ws.Cells.Item(43, 52).Value = 0
$ws.Cells.Item(43, 53).Value = 0
$ws.Cells.Item(43, 54).Value = 0
$ws.Cells.Item(43, 55).Value = 0
$ws.Cells.Item(43, 56).Value = 44.21
$ws.Cells.Item(43, 57).Value = 321.19
$ws.Cells.Item(43, 58).Value = 0.61
$ws.Cells.Item(43, 59).Value = 0

# Row 44
$ws.Cells.Item(44, 1).Value = 66
$ws.Cells.Item(44, 2).Value = 0
$ws.Cells.Item(44, 3).Value = 1.6
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 1.59
$ws.Cells.Item(44, 9).Value = 13.01
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(44, 15).Value = 12.96
$ws.Cells.Item(44, 16).Value = 12.96
$ws.Cells.Item(44, 17).Value = 12.96
$ws.Cells.Item(44, 18).Value = 0
$ws.Cells.Item(44, 19).Value = 0
$ws.Cells.Item(44, 20).Value = 0
$ws.Cells.Item(44, 21).Value = 12.96
$ws.Cells.Item(44, 22).Value = 0
$ws.Cells.Item(44, 23).Value = 0
$ws.Cells.Item(44, 24).Value = 0
$ws.Cells.Item(44, 25).Value = 0
$ws.Cells.Item(44, 26).Value = 0
$ws.Cells.Item(44, 27).Value = 49.31
$ws.Cells.Item(44, 28).Value = 321.52
$ws.Cells.Item(44, 29).Value = 3.24
$ws.Cells.Item(44, 30).Value = 0
$ws.Cells.Item(44, 31).Value = 0
$ws.Cells.Item(44, 32).Value = 1.59
$ws.Cells.Item(44, 33).Value = 0
$ws.Cells.Item(44, 34).Value = 0
$ws.Cells.Item(44, 35).Value = 0
$ws.Cells.Item(44, 36).Value = 0
$ws.Cells.Item(44, 37).Value = 1.56
$ws.Cells.Item(44, 38).Value = 13.03
$ws.Cells.Item(44, 39).Value = 0
$ws.Cells.Item(44, 40).Value = 0
$ws.Cells.Item(44, 41).Value = 0
$ws.Cells.Item(44, 42).Value = 0
$ws.Cells.Item(44, 43).Value = 0
$ws.Cells.Item(44, 44).Value = 12.99
$ws.Cells.Item(44, 45).Value = 12.99
$ws.Cells.Item(44, 46).Value = 12.99
$ws.Cells.Item(44, 47).Value = 0
$ws.Cells.Item(44, 48).Value = 0
$ws.Cells.Item(44, 49).Value = 0
$ws.Cells.Item(44, 50).Value = 12.99
$ws.Cells.Item(44, 51).Value = 0
$ws.Cells.Item(44, 52).Value = 0
$ws.Cells.Item(44, 53).Value = 0
$ws.Cells.Item(44, 54).Value = 0
$ws.Cells.Item(44, 55).Value = 0
$ws.Cells.Item(44, 56).Value = 47.37
$ws.Cells.Item(44, 57).Value = 323.46
$ws.Cells.Item(44, 58).Value = 3.55
$ws.Cells.Item(44, 59).Value = 0

# Row 45
$ws.Cells.Item(45, 1).Value = 68
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = 4.76
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 4.76
$ws.Cells.Item(45, 9).Value = 13
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 13
$ws.Cells.Item(45, 16).Value = 13
$ws.Cells.Item(45, 17).Value = 13
$ws.Cells.Item(45, 18).Value = 0
$ws.Cells.Item(45, 19).Value = 0
$ws.Cells.Item(45, 20).Value = 0
$ws.Cells.Item(45, 21).Value = 13
$ws.Cells.Item(45, 22).Value = 0
$ws.Cells.Item(45, 23).Value = 0
$ws.Cells.Item(45, 24).Value = 0
$ws.Cells.Item(45, 25).Value = 0
$ws.Cells.Item(45, 26).Value = 0
$ws.Cells.Item(45, 27).Value = 69.40000000000001
$ws.Cells.Item(45, 28).Value = 295.02
$ws.Cells.Item(45, 29).Value = 1.32
$ws.Cells.Item(45, 30).Value = 0
$ws.Cells.Item(45, 31).Value = 0
$ws.Cells.Item(45, 32).Value = 4.78
$ws.Cells.Item(45, 33).Value = 0
$ws.Cells.Item(45, 34).Value = 0
$ws.Cells.Item(45, 35).Value = 0
$ws.Cells.Item(45, 36).Value = 0
$ws.Cells.Item(45, 37).Value = 4.77
$ws.Cells.Item(45, 38).Value = 13.01
$ws.Cells.Item(45, 39).Value = 0
$ws.Cells.Item(45, 40).Value = 0
$ws.Cells.Item(45, 41).Value = 0
$ws.Cells.Item(45, 42).Value = 0
$ws.Cells.Item(45, 43).Value = 0
$ws.Cells.Item(45, 44).Value = 13
$ws.Cells.Item(45, 45).Value = 13
$ws.Cells.Item(45, 46).Value = 13
$ws.Cells.Item(45, 47).Value = 0
$ws.Cells.Item(45, 48).Value = 0
$ws.Cells.Item(45, 49).Value = 0
$ws.Cells.Item(45, 50).Value = 13
$ws.Cells.Item(45, 51).Value = 0
$ws.Cells.Item(45, 52).Value = 0
$ws.Cells.Item(45, 53).Value = 0
$ws.Cells.Item(45, 54).Value = 0
$ws.Cells.Item(45, 55).Value = 0
$ws.Cells.Item(45, 56).Value = 68.95
$ws.Cells.Item(45, 57).Value = 295.45
$ws.Cells.Item(45, 58).Value = 1.48
$ws.Cells.Item(45, 59).Value = 0

# Row 46
$ws.Cells.Item(46, 1).Value = 70
$ws.Cells.Item(46, 2).Value = 0.01
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 0.01
$ws.Cells.Item(46, 7).Value = 16
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 16
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = 16
$ws.Cells.Item(46, 14).Value = 16
$ws.Cells.Item(46, 15).Value = 16.93
$ws.Cells.Item(46, 16).Value = 16.93
$ws.Cells.Item(46, 17).Value = 16.93
$ws.Cells.Item(46, 18).Value = 0
$ws.Cells.Item(46, 19).Value = 0
$ws.Cells.Item(46, 20).Value = 16.93
$ws.Cells.Item(46, 21).Value = 0
$ws.Cells.Item(46, 22).Value = 16.87
$ws.Cells.Item(46, 23).Value = 0
$ws.Cells.Item(46, 24).Value = 0
$ws.Cells.Item(46, 25).Value = 0
$ws.Cells.Item(46, 26).Value = 0
$ws.Cells.Item(46, 27).Value = 1
$ws.Cells.Item(46, 28).Value = 300.18
$ws.Cells.Item(46, 29).Value = 0.5600000000000001
$ws.Cells.Item(46, 30).Value = 0
$ws.Cells.Item(46, 31).Value = 0
$ws.Cells.Item(46, 32).Value = 0
$ws.Cells.Item(46, 33).Value = 0
$ws.Cells.Item(46, 34).Value = 16
$ws.Cells.Item(46, 35).Value = 0
$ws.Cells.Item(46, 36).Value = 17
$ws.Cells.Item(46, 37).Value = 0
$ws.Cells.Item(46, 38).Value = 0
$ws.Cells.Item(46, 39).Value = 0
$ws.Cells.Item(46, 40).Value = 16
$ws.Cells.Item(46, 41).Value = 0
$ws.Cells.Item(46, 42).Value = 16
$ws.Cells.Item(46, 43).Value = 16
$ws.Cells.Item(46, 44).Value = 16.95
$ws.Cells.Item(46, 45).Value = 16.95
$ws.Cells.Item(46, 46).Value = 16.95
$ws.Cells.Item(46, 47).Value = 0
$ws.Cells.Item(46, 48).Value = 0
$ws.Cells.Item(46, 49).Value = 16.95
$ws.Cells.Item(46, 50).Value = 0
$ws.Cells.Item(46, 51).Value = 16
$ws.Cells.Item(46, 52).Value = 0
$ws.Cells.Item(46, 53).Value = 0
$ws.Cells.Item(46, 54).Value = 0
$ws.Cells.Item(46, 55).Value = 0
$ws.Cells.Item(46, 56).Value = 0.98
$ws.Cells.Item(46, 57).Value = 301.07
$ws.Cells.Item(46, 58).Value = 0.66
$ws.Cells.Item(46, 59).Value = 0

# Row 47
$ws.Cells.Item(47, 1).Value = 71
$ws.Cells.Item(47, 2).Value = 0.02
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 16.01
$ws.Cells.Item(47, 6).Value = 0.02
$ws.Cells.Item(47, 7).Value = 17
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 16
$ws.Cells.Item(47, 12).Value = 0.01
$ws.Cells.Item(47, 13).Value = 16
$ws.Cells.Item(47, 14).Value = 16
$ws.Cells.Item(47, 15).Value = 16.96
$ws.Cells.Item(47, 16).Value = 16.96
$ws.Cells.Item(47, 17).Value = 16.96
$ws.Cells.Item(47, 18).Value = 0
$ws.Cells.Item(47, 19).Value = 0
$ws.Cells.Item(47, 20).Value = 16.96
$ws.Cells.Item(47, 21).Value = 0
$ws.Cells.Item(47, 22).Value = 16.86
$ws.Cells.Item(47, 23).Value = 0
$ws.Cells.Item(47, 24).Value = 0
$ws.Cells.Item(47, 25).Value = 0
$ws.Cells.Item(47, 26).Value = 0
$ws.Cells.Item(47, 27).Value = 2.26
$ws.Cells.Item(47, 28).Value = 298.85
$ws.Cells.Item(47, 29).Value = 0.06
$ws.Cells.Item(47, 30).Value = 0
$ws.Cells.Item(47, 31).Value = 0.06
$ws.Cells.Item(47, 32).Value = 0
$ws.Cells.Item(47, 33).Value = 0
$ws.Cells.Item(47, 34).Value = 15.99
$ws.Cells.Item(47, 35).Value = 0.06
$ws.Cells.Item(47, 36).Value = 16.97
$ws.Cells.Item(47, 37).Value = 0
$ws.Cells.Item(47, 38).Value = 0
$ws.Cells.Item(47, 39).Value = 0
$ws.Cells.Item(47, 40).Value = 16
$ws.Cells.Item(47, 41).Value = 0
$ws.Cells.Item(47, 42).Value = 15.99
$ws.Cells.Item(47, 43).Value = 15.99
$ws.Cells.Item(47, 44).Value = 16.57
$ws.Cells.Item(47, 45).Value = 16.57
$ws.Cells.Item(47, 46).Value = 16.57
$ws.Cells.Item(47, 47).Value = 0
$ws.Cells.Item(47, 48).Value = 0
$ws.Cells.Item(47, 49).Value = 16.57
$ws.Cells.Item(47, 50).Value = 0
$ws.Cells.Item(47, 51).Value = 15.99
$ws.Cells.Item(47, 52).Value = 0
$ws.Cells.Item(47, 53).Value = 0
$ws.Cells.Item(47, 54).Value = 0
$ws.Cells.Item(47, 55).Value = 0
$ws.Cells.Item(47, 56).Value = 2.44
$ws.Cells.Item(47, 57).Value = 299.73
$ws.Cells.Item(47, 58).Value = 0.02
$ws.Cells.Item(47, 59).Value = 0
